# Invitee sample: add new "Profile picture" column header (import/export profile pic url)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cell P1 - "Profile picture"
$ws.Range("P1").Value = "Profile picture"

# Keep selection consistent with the new last-used cell (matches source workbook behavior)
$ws.Range("P2").Select()
